$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.998.58"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "2.420.97"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.54"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.50"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.109"
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("E11").Value = "  -3.41%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.82"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "2.858.36"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "61.951.13"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "2.420.05"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.27"
$ws.Range("E18").Value = "  +1.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "323.63"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.86"
$ws.Range("E20").Value = "  +2.63%  "
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.58"
$ws.Range("E23").Value = "  +2.50%  "
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.75"
$ws.Range("E25").Value = "  -2.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "563.62"
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").Value = "2.540.88"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "0.0₃0943"
$ws.Range("E29").Value = "  +2.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.23"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("E31").Value = "  -3.16%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("E34").Value = "  -3.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.79"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "154.32"
$ws.Range("E38").Value = "  +2.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.46"
$ws.Range("E39").Value = "  -3.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.55"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "148.23"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.64"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.87"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.594"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.56"
$ws.Range("E51").Value = "  +0.76%  "
